$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Initial Investment ---
$ws.Range("B2").Value = -2500000

# --- Row 3: Depreciation (C3:L3) ---
$ws.Range("C3:L3").Value = 75000

# --- Row 4: Incoming Payments (C4:L4) ---
$ws.Range("C4:L4").Value = 1400000

# --- Row 5: Outgoing Payments ---
$ws.Range("B5").Value = -1050000
$ws.Range("C5:L5").Value = -700000

# --- Row 6: Residual ---
$ws.Range("L6").Value = 140000

# --- Row 7: restricted Equity ---
$ws.Range("B7").Value = -200000
$ws.Range("L7").Value = 200000

# --- Row 8: Yearly Net ---
$ws.Range("B8").Value = -3750000
$ws.Range("C8:K8").Value = 775000
$ws.Range("L8").Value = 1115000

# --- Row 9: Present Value ---
$ws.Range("B9").Value = -3750000
$ws.Range("C9").Value = 714944.6494464944
$ws.Range("D9").Value = 659543.0345447365
$ws.Range("E9").Value = 608434.5337128565
$ws.Range("F9").Value = 561286.4702148122
$ws.Range("G9").Value = 517791.9466926312
$ws.Range("H9").Value = 477667.8475024272
$ws.Range("I9").Value = 440652.9958509476
$ws.Range("J9").Value = 406506.4537370366
$ws.Range("K9").Value = 375005.9536319525
$ws.Range("L9").Value = 497716.5079152804

# --- Row 10: Accumulated Present Value ---
$ws.Range("B10").Value = -3750000
$ws.Range("C10").Value = -3035055.350553506
$ws.Range("D10").Value = -2375512.316008769
$ws.Range("E10").Value = -1767077.782295913
$ws.Range("F10").Value = -1205791.3120811
$ws.Range("G10").Value = -687999.3653884693
$ws.Range("H10").Value = -210331.5178860421
$ws.Range("I10").Value = 230321.4779649056
$ws.Range("J10").Value = 636827.9317019421
$ws.Range("K10").Value = 1011833.885333895
$ws.Range("L10").Value = 1509550.393249175

# F10:H10 turned negative, so they need to switch from the "positive"
# (green fill) style to the "negative" (red fill) style, matching E10's
# formatting. Copy/PasteSpecial of formats reuses the existing style
# definition instead of fabricating a brand-new one.
$ws.Range("E10").Copy()
$ws.Range("F10:H10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 11: Net Present Value ---
$ws.Range("B11").Value = 1509550.393249175
